# Auto-generated edit script: applies the Seraph_Profits market-data refresh diff
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets (the xlsx sheet tabs that
# correspond one-to-one, in tab order, with the commit's "Sheets/Seraph_Profits.xlsx"
# rows). For each touched cell we either set its new numeric value, or (where the
# diff shows the cell disappearing entirely) clear it so it is dropped from the XML.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1004.4545
$ws.Range("I2").Value = 987.1429000000001
$ws.Range("J2").Value = 1034.75
$ws.Range("K2").Value = 987.1429000000001
$ws.Range("L2").Value = 1034.75
$ws.Range("M2").Value = -874.1429000000001
$ws.Range("N2").Value = -1260.75
$ws.Range("H17").Value = 1249.8182
$ws.Range("J17").Value = 1258.7142
$ws.Range("L17").Value = 3776.1426
$ws.Range("N17").Value = -4112.142599999999
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H41").Value = 176.72728
$ws.Range("I41").Value = 176.72728
$ws.Range("K41").Value = 176.72728
$ws.Range("M41").Value = 263.27272
$ws.Range("H52").Value = 4999.5
$ws.Range("J52").Value = 4999.5
$ws.Range("L52").Value = 14998.5
$ws.Range("N52").Value = -15318.5
$ws.Range("H112").Value = 2267.111
$ws.Range("J112").Value = 2607.7144
$ws.Range("L112").Value = 7823.1432
$ws.Range("N112").Value = -10039.1432
$ws.Range("H125").Value = 1008.8
$ws.Range("I125").Value = 1021.3333
$ws.Range("J125").Value = 990
$ws.Range("K125").Value = 9191.9997
$ws.Range("L125").Value = 8910
$ws.Range("M125").Value = -6731.9997
$ws.Range("N125").Value = -13830
$ws.Range("H137").Value = 1635.0714
$ws.Range("I137").Value = 1486.625
$ws.Range("K137").Value = 4459.875
$ws.Range("M137").Value = -1909.875
$ws.Range("H138").Value = 3763.4146
$ws.Range("I138").Value = 3070.3
$ws.Range("K138").Value = 9210.900000000001
$ws.Range("M138").Value = -4070.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H74").Value = 1199.3334
$ws.Range("I74").Value = 1199.3334
$ws.Range("K74").Value = 1199.3334
$ws.Range("M74").Value = -325.3334
$ws.Range("H77").Value = 1199.3334
$ws.Range("I77").Value = 1199.3334
$ws.Range("K77").Value = 5996.666999999999
$ws.Range("M77").Value = -1628.666999999999
$ws.Range("H122").Value = 2614.25
$ws.Range("I122").Value = 2069.3333
$ws.Range("J122").Value = 4249
$ws.Range("K122").Value = 6207.999899999999
$ws.Range("L122").Value = 12747
$ws.Range("M122").Value = -3757.999899999999
$ws.Range("N122").Value = -17647
$ws.Range("H132").Value = 1493.3636
$ws.Range("I132").Value = 1493
$ws.Range("K132").Value = 4479
$ws.Range("M132").Value = -1949
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1718.7142
$ws.Range("I20").Value = 1547
$ws.Range("K20").Value = 1547
$ws.Range("M20").Value = -1300
$ws.Range("H97").Value = 36299.668
$ws.Range("I97").Value = 22450
$ws.Range("K97").Value = 22450
$ws.Range("M97").Value = -21459
$ws.Range("H105").Value = 2556.2222
$ws.Range("I105").Value = 2571.9167
$ws.Range("J105").Value = 2524.8333
$ws.Range("K105").Value = 2571.9167
$ws.Range("L105").Value = 2524.8333
$ws.Range("M105").Value = -824.9167000000002
$ws.Range("N105").Value = -6018.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 242.5
$ws.Range("I22").Value = 242.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 242.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 107.5
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 2900.5
$ws.Range("I31").Value = 2115.8572
$ws.Range("J31").Value = 5254.4287
$ws.Range("K31").Value = 2115.8572
$ws.Range("L31").Value = 5254.4287
$ws.Range("M31").Value = -1820.8572
$ws.Range("N31").Value = -5844.4287
$ws.Range("H34").Value = 2900.5
$ws.Range("I34").Value = 2115.8572
$ws.Range("J34").Value = 5254.4287
$ws.Range("K34").Value = 2115.8572
$ws.Range("L34").Value = 5254.4287
$ws.Range("M34").Value = -1913.8572
$ws.Range("N34").Value = -5658.4287
$ws.Range("H58").Value = 2073.923
$ws.Range("J58").Value = 4491
$ws.Range("L58").Value = 4491
$ws.Range("N58").Value = -4897
$ws.Range("H70").Value = 38332.668
$ws.Range("J70").Value = 38332.668
$ws.Range("L70").Value = 38332.668
$ws.Range("N70").Value = -38962.668
$ws.Range("H73").Value = 38332.668
$ws.Range("J73").Value = 38332.668
$ws.Range("L73").Value = 38332.668
$ws.Range("N73").Value = -40516.668
$ws.Range("H97").Value = 63999
$ws.Range("J97").Value = 63999
$ws.Range("L97").Value = 63999
$ws.Range("N97").Value = -65981
$ws.Range("H107").Value = 831.5714
$ws.Range("I107").Value = 853.6667
$ws.Range("J107").Value = 699
$ws.Range("K107").Value = 853.6667
$ws.Range("L107").Value = 699
$ws.Range("M107").Value = 1066.3333
$ws.Range("N107").Value = -4539
$ws.Range("H132").Value = 1315.8
$ws.Range("I132").Value = 1393
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 4179
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -1649
$ws.Range("N132").Value = -8660
$ws.Range("H134").Value = 2137.12
$ws.Range("I134").Value = 1830.8572
$ws.Range("K134").Value = 5492.571599999999
$ws.Range("M134").Value = -2957.571599999999
$ws.Range("H136").Value = 2073.923
$ws.Range("J136").Value = 4491
$ws.Range("L136").Value = 13473
$ws.Range("N136").Value = -18573

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1425.5834
$ws.Range("I5").Value = 781.1429000000001
$ws.Range("K5").Value = 2343.4287
$ws.Range("M5").Value = -2231.4287
$ws.Range("H92").Value = 286.94446
$ws.Range("J92").Value = 303.57144
$ws.Range("L92").Value = 910.71432
$ws.Range("N92").Value = -3406.71432
$ws.Range("H104").Value = 107999.4
$ws.Range("J104").Value = 107999.4
$ws.Range("L104").Value = 323998.2
$ws.Range("N104").Value = -329240.2
$ws.Range("H131").Value = 897.5
$ws.Range("I131").Value = 342.5
$ws.Range("K131").Value = 1027.5
$ws.Range("M131").Value = 4012.5
$ws.Range("H135").Value = 1425.5834
$ws.Range("I135").Value = 781.1429000000001
$ws.Range("K135").Value = 7030.2861
$ws.Range("M135").Value = -4495.2861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7374.8125
$ws.Range("I70").Value = 6666.5557
$ws.Range("K70").Value = 6666.5557
$ws.Range("M70").Value = -6396.5557
$ws.Range("H73").Value = 7374.8125
$ws.Range("I73").Value = 6666.5557
$ws.Range("K73").Value = 6666.5557
$ws.Range("M73").Value = -5730.5557
$ws.Range("H92").Value = 14363.4
$ws.Range("J92").Value = 11704.25
$ws.Range("L92").Value = 11704.25
$ws.Range("N92").Value = -15448.25
$ws.Range("H102").Value = 2561.125
$ws.Range("I102").Value = 2512.7144
$ws.Range("K102").Value = 2512.7144
$ws.Range("M102").Value = -890.7143999999998
$ws.Range("H109").Value = 59500
$ws.Range("J109").Value = 59500
$ws.Range("L109").Value = 59500
$ws.Range("N109").Value = -61580
$ws.Range("H122").Value = 35286.332
$ws.Range("I122").Value = 1317.625
$ws.Range("K122").Value = 3952.875
$ws.Range("M122").Value = -1502.875
$ws.Range("H132").Value = 2697.3333
$ws.Range("I132").Value = 2325.1428
$ws.Range("K132").Value = 6975.428400000001
$ws.Range("M132").Value = -4445.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2280.0908
$ws.Range("I40").Value = 2310.1
$ws.Range("J40").Value = 1980
$ws.Range("K40").Value = 2310.1
$ws.Range("L40").Value = 1980
$ws.Range("M40").Value = -2174.1
$ws.Range("N40").Value = -2252
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6968.625
$ws.Range("I62").Value = 4125
$ws.Range("K62").Value = 4125
$ws.Range("M62").Value = -3501
$ws.Range("H65").Value = 6968.625
$ws.Range("I65").Value = 4125
$ws.Range("K65").Value = 20625
$ws.Range("M65").Value = -17505
$ws.Range("H96").Value = 3000
$ws.Range("I96").Value = 3000
$ws.Range("K96").Value = 3000
$ws.Range("M96").Value = -1627
